# ============================================================
# Edit: add "Sheet3" with two Snakes parameter tables, update
# workbook active tab and selections (per commit:
# "working on different parameters for snakes")
# ============================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- 1. Adjust Sheet1 (was the active/selected tab) ---------------
$ws1.Select()
$ws1.Range("A1").Select()

# --- 2. Adjust Sheet2 selection -----------------------------------
$ws2.Select()
$ws2.Range("J15").Select()

# --- 3. Add new Sheet3 right after Sheet2 -------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"

# Colours (VBA BGR-packed ints) reused from the existing red font:
#   RGB(FF,33,00) -> 0x0033FF -> 13311
$redColor = 13311

# style codes used below (mirrors the workbook's existing cellXfs):
#   0 = General number format, default font
#   1 = General number format, bold font            (section headers)
#   2 = "#,##0.00" number format, default font       (plain decimals)
#   3 = "0.00%" number format, red font              (percentages)
#   4 = "0.00%" number format, red BOLD font         (percentages, emphasised)

function Set-CellStyle($cell, [string]$style) {
    switch ($style) {
        "0" {
            $cell.NumberFormat = "General"
        }
        "1" {
            $cell.NumberFormat = "General"
            $cell.Font.Bold = $true
        }
        "2" {
            $cell.NumberFormat = "#,##0.00"
        }
        "3" {
            $cell.NumberFormat = "0.00%"
            $cell.Font.Color = $redColor
        }
        "4" {
            $cell.NumberFormat = "0.00%"
            $cell.Font.Color = $redColor
            $cell.Font.Bold = $true
        }
    }
}

# --- 4. Populate Sheet3 with the two parameter tables -------------
# Pipe-delimited: Row|Col|StyleCode|Kind|Payload
#   Kind S = shared text, N = literal number, F = formula (A1-style, no leading '=')
$cellData = @"
1|A|1|S|Snakes – 2412-01
1|B|0|S|lung case 1
1|C|0|S|lung case 2
1|D|0|S|lung case 3
1|E|0|S|lung case 4
1|F|0|S|lung case 5
1|G|0|S|lung average
1|H|0|S|lung max
1|I|0|S|lung min
2|A|0|S|volume of the actual possible
2|B|2|N|52295.681
2|C|2|N|25264.893
2|D|2|N|22822.453
2|E|2|N|25799.604
2|F|2|N|50882.94
2|G|0|F|AVERAGE(B2:F2)
2|H|0|F|MAX(B2:F2)
2|I|0|F|MIN(B2:F2)
3|A|0|S|volume of the estimated possible
3|B|2|N|49914.849
3|C|2|N|23469.238
3|D|2|N|20752.461
3|E|2|N|24804.792
3|F|2|N|47637.539
3|G|0|F|AVERAGE(B3:F3)
3|H|0|F|MAX(B3:F3)
3|I|0|F|MIN(B3:F3)
4|A|3|S|percentage volume error of the actual and estimated possible
4|B|3|N|0.046
4|C|3|N|0.071
4|D|3|N|0.091
4|E|3|N|0.039
4|F|3|N|0.064
4|G|3|F|AVERAGE(B4:F4)
4|H|3|F|MAX(B4:F4)
4|I|3|F|MIN(B4:F4)
5|A|0|S|volume of the actual consensus
5|B|0|N|32805.908
5|C|0|N|9591.064
5|D|0|N|11771.3
5|E|0|N|15589.534
5|F|0|N|34140.034
5|G|0|F|AVERAGE(B5:F5)
5|H|0|F|MAX(B5:F5)
5|I|0|F|MIN(B5:F5)
6|A|0|S|volume of the estimated consensus
6|B|2|N|35978.927
6|C|2|N|15380.249
6|D|2|N|14265.971
6|E|2|N|19350.482
6|F|2|N|40749.531
6|G|0|F|AVERAGE(B6:F6)
6|H|0|F|MAX(B6:F6)
6|I|0|F|MIN(B6:F6)
7|A|3|S|percentage volume error of the actual and estimated consensus
7|B|3|N|-0.097
7|C|4|N|-0.604
7|D|4|N|-0.212
7|E|4|F|-0.241
7|F|4|N|-0.194
7|G|4|F|AVERAGE(B7:F7)
7|H|3|F|MAX(B7:F7)
7|I|3|F|MIN(B7:F7)
8|A|0|S|volume of the actual variability
8|B|2|N|19489.773
8|C|2|N|15673.828
8|D|2|N|11051.153
8|E|2|N|10210.069
8|F|2|N|16742.907
8|G|0|F|AVERAGE(B8:F8)
8|H|0|F|MAX(B8:F8)
8|I|0|F|MIN(B8:F8)
9|A|0|S|volume of the estimated variability
9|B|2|N|13935.921
9|C|2|N|8088.989
9|D|2|N|6486.49
9|E|2|N|5454.31
9|F|2|N|6888.008
9|G|0|F|AVERAGE(B9:F9)
9|H|0|F|MAX(B9:F9)
9|I|0|F|MIN(B9:F9)
10|A|3|S|percentage volume error of the actual and estimated variability
10|B|4|N|0.285
10|C|4|N|0.484
10|D|4|N|0.413
10|E|4|N|0.466
10|F|4|N|0.589
10|G|4|F|AVERAGE(B10:F10)
10|H|4|F|MAX(B10:F10)
10|I|4|F|MIN(B10:F10)
13|A|1|S|Snakes – 2412-02
13|B|0|S|lung case 1
13|C|0|S|lung case 2
13|D|0|S|lung case 3
13|E|0|S|lung case 4
13|F|0|S|lung case 5
13|G|0|S|lung average
13|H|0|S|lung max
13|I|0|S|lung min
14|A|0|S|volume of the actual possible
14|B|2|N|52295.681
14|C|2|N|25264.893
14|D|2|N|22822.453
14|E|2|N|25799.604
14|F|2|N|50882.94
14|G|0|F|AVERAGE(B14:F14)
14|H|0|F|MAX(B14:F14)
14|I|0|F|MIN(B14:F14)
15|A|0|S|volume of the estimated possible
15|B|2|N|49927.672
15|C|2|N|23166.504
15|D|2|N|20624.11
15|E|2|N|24829.741
15|F|2|N|47423.983
15|G|0|F|AVERAGE(B15:F15)
15|H|0|F|MAX(B15:F15)
15|I|0|F|MIN(B15:F15)
16|A|3|S|percentage volume error of the actual and estimated possible
16|B|3|N|0.045
16|C|3|N|0.083
16|D|3|N|0.096
16|E|3|N|0.038
16|F|3|N|0.068
16|G|3|F|AVERAGE(B16:F16)
16|H|3|F|MAX(B16:F16)
16|I|3|F|MIN(B16:F16)
17|A|0|S|volume of the actual consensus
17|B|0|N|32805.908
17|C|0|N|9591.064
17|D|0|N|11771.3
17|E|0|N|15589.534
17|F|0|N|34140.034
17|G|0|F|AVERAGE(B17:F17)
17|H|0|F|MAX(B17:F17)
17|I|0|F|MIN(B17:F17)
18|A|0|S|volume of the estimated consensus
18|B|2|N|29906.45
18|C|2|N|13271.484
18|D|2|N|11724.783
18|E|2|N|16094.736
18|F|2|N|37310.2
18|G|0|F|AVERAGE(B18:F18)
18|H|0|F|MAX(B18:F18)
18|I|0|F|MIN(B18:F18)
19|A|3|S|percentage volume error of the actual and estimated consensus
19|B|3|N|0.088
19|C|4|N|-0.384
19|D|3|N|0.004
19|E|3|F|-0.032
19|F|3|N|-0.093
19|G|3|F|AVERAGE(B19:F19)
19|H|3|F|MAX(B19:F19)
19|I|3|F|MIN(B19:F19)
20|A|0|S|volume of the actual variability
20|B|2|N|19489.773
20|C|2|N|15673.828
20|D|2|N|11051.153
20|E|2|N|10210.069
20|F|2|N|16742.907
20|G|0|F|AVERAGE(B20:F20)
20|H|0|F|MAX(B20:F20)
20|I|0|F|MIN(B20:F20)
21|A|0|S|volume of the estimated variability
21|B|2|N|20021.221
21|C|2|N|9895.02
21|D|2|N|8899.327
21|E|2|N|8735.004
21|F|2|N|10113.783
21|G|0|F|AVERAGE(B21:F21)
21|H|0|F|MAX(B21:F21)
21|I|0|F|MIN(B21:F21)
22|A|3|S|percentage volume error of the actual and estimated variability
22|B|3|N|-0.027
22|C|4|N|0.369
22|D|4|N|0.195
22|E|4|N|0.144
22|F|4|N|0.396
22|G|4|F|AVERAGE(B22:F22)
22|H|4|F|MAX(B22:F22)
22|I|3|F|MIN(B22:F22)
"@

$rows = $cellData -split "`n"
foreach ($row in $rows) {
    $row = $row.Trim()
    if ($row.Length -eq 0) { continue }
    $parts = $row -split '\|', 5
    $r = [int]$parts[0]
    $col = $parts[1]
    $style = $parts[2]
    $kind = $parts[3]
    $payload = $parts[4]

    $addr = "$col$r"
    $cell = $ws3.Range($addr)

    if ($kind -eq "F") {
        $cell.Formula = "=" + $payload
    } elseif ($kind -eq "N") {
        $cell.Value = [double]$payload
    } else {
        $cell.Value = $payload
    }

    Set-CellStyle $cell $style
}

# --- 5. Sheet3 view / selection / column widths --------------------
$ws3.Range("A1:A22").ColumnWidth = 52.9897959183674
$ws3.Range("B1:I22").ColumnWidth = 11.5204081632653

$ws3.Activate()
$ws3.Range("I22").Select()

# --- 6. Re-apply Sheet1 / Sheet2 column widths (tiny re-measure) --
$ws1.Range("A1:A46").ColumnWidth = 53.2959183673469
$ws1.Range("B1:U46").ColumnWidth = 11.7551020408163

$ws2.Range("A1:A11").ColumnWidth = 54.3826530612245
$ws2.Range("B1:AK11").ColumnWidth = 11.7551020408163
